$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row number + the new values for any changed columns (B/C/D/E).
# D (price) values are text that can look numeric ("1.008", "3.200", "9.950", ...);
# forcing NumberFormat "@" before the write and resetting the Style afterwards keeps
# them as literal text (no precision loss, e.g. "3.200" staying "3.200" not 3.2) while
# not leaving a lasting number-format style on the cell.
$rows = @(
    @{Row=2; D="25.890.70"; E="  +0.11%  "},
    @{Row=3; D="1.639.80"; E="  +0.44%  "},
    @{Row=4; D="1.008"; E="  -0.12%  "},
    @{Row=5; D="214.92"; E="  -0.13%  "},
    @{Row=6; D="0.5035"; E="  +0.59%  "},
    @{Row=7; D="1.005"; E="  -0.53%  "},
    @{Row=8; D="0.2574"; E="  +0.23%  "},
    @{Row=9; D="0.06386"; E="  -0.39%  "},
    @{Row=10; D="19.56"; E="  +0.52%  "},
    @{Row=11; D="0.07797"; E="  +0.83%  "},
    @{Row=12; D="1.656.72"; E="  +1.44%  "},
    @{Row=13; D="4.275"; E="  +0.69%  "},
    @{Row=14; D="1.862.97"; E="  +0.22%  "},
    @{Row=15; D="0.5428"; E="  -0.12%  "},
    @{Row=16; D="0.0₅7877"; E="  -0.43%  "},
    @{Row=17; D="64.82"; E="  +2.15%  "},
    @{Row=18; D="25.944.00"; E="  +0.17%  "},
    @{Row=19; D="1.007"; E="  -0.30%  "},
    @{Row=20; D="197.43"; E="  -2.93%  "},
    @{Row=21; D="4.394"; E="  +2.41%  "},
    @{Row=22; D="9.950"; E="  -0.31%  "},
    @{Row=23; D="5.976"; E="  +0.76%  "},
    @{Row=24; D="1.006"; E="  -0.51%  "},
    @{Row=25; D="1.882"; E="  -3.91%  "},
    @{Row=26; D="140.01"; E="  -0.62%  "},
    @{Row=27; D="0.1143"; E="  -0.53%  "},
    @{Row=28; D="6.844"; E="  +1.18%  "},
    @{Row=29; D="15.69"; E="  -0.35%  "},
    @{Row=30; D="1.241"; E="  +0.14%  "},
    @{Row=31; D="0.04872"; E="  -4.19%  "},
    @{Row=32; D="3.261"; E="  +0.12%  "},
    @{Row=33; D="3.200"; E="  +0.41%  "},
    @{Row=34; D="1.536"; E="  -0.27%  "},
    @{Row=35; D="2.373"; E="  +1.29%  "},
    @{Row=36; D="0.8896"; E="  -0.07%  "},
    @{Row=37; D="2.612"; E="  +0.23%  "},
    @{Row=38; D="0.5539"; E="  -1.86%  "},
    @{Row=39; D="1.132.45"; E="  -0.27%  "},
    @{Row=40; D="0.01559"; E="  +0.14%  "},
    @{Row=41; D="1.005"; E="  -0.50%  "},
    @{Row=42; D="5.695"; E="  +1.05%  "},
    @{Row=43; D="0.8156"; E="  -0.02%  "},
    @{Row=44; D="99.55"; E="  +0.10%  "},
    @{Row=45; B="RocketPoolETH"; C="https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"; D="1.773.30"; E="  +0.20%  "},
    @{Row=46; B="BabyDogeCoin"; C="https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"; D="0.0₈118"; E="  +5.01%  "},
    @{Row=47; D="0.4537"; E="  +0.27%  "},
    @{Row=48; B="Aave"; C="https://coinranking.com/coin/ixgUfzmLR+aave-aave"; D="55.37"; E="  +1.34%  "},
    @{Row=49; B="Frax"; C="https://coinranking.com/coin/KfWtaeV1W+frax-frax"; D="1.008"; E="  -0.34%  "},
    @{Row=50; D="0.05089"; E="  +1.38%  "},
    @{Row=51; D="1.007"; E="  -0.24%  "}
)

function Set-TextCell($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

foreach ($r in $rows) {
    if ($r.ContainsKey("B")) { Set-TextCell $ws.Cells.Item($r.Row, 2) $r.B }
    if ($r.ContainsKey("C")) { Set-TextCell $ws.Cells.Item($r.Row, 3) $r.C }
    if ($r.ContainsKey("D")) { Set-TextCell $ws.Cells.Item($r.Row, 4) $r.D }
    if ($r.ContainsKey("E")) { Set-TextCell $ws.Cells.Item($r.Row, 5) $r.E }
}
